$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the precondition text for the housing_transform row (C4)
$ws.Range("C4").Value = "?ARk <= ?C(?Rk)"

# Update the active cell selection to B14
$ws.Range("B14").Select()
